# bugfix: 修复不生成 consts.json 的bug
#
# Adds a new "nested" field (column R) of type `array:map:int,string` to the
# test schema sheet, and fixes the sample value in L7 from the English
# placeholder to the Chinese one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column R: array:map:int,string / nested / cs / 嵌套类型 ---
$ws.Range("R1").Value = "array:map:int,string"
$ws.Range("R2").Value = "nested"
$ws.Range("R3").Value = "cs"
$ws.Range("R4").Value = "嵌套类型"
$ws.Range("R6").Value = '[{1:"Happy",2:"Smile"}]'
$ws.Range("R7").Value = '[{-1:"开心",-2:"笑"}]'

# --- Fix sample string value in L7 ---
$ws.Range("L7").Value = '"你好 Alpaca"'

# --- Column width for the new column R ---
$ws.Columns.Item(18).ColumnWidth = 22.57

# --- Restore the active selection cell shown in the saved view ---
$ws.Range("O11").Select() | Out-Null
